# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# The "Periodo Mora" column (E16:E32) is re-ordered from descending
# (2401 .. 2209) to ascending (2209 .. 2401).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2209","2210","2211","2212","2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312","2401")

$row = 16
foreach ($p in $periods) {
    $ws.Range("E$row").Value = $p
    $row++
}
